# Insert a new row above row 4 on the first worksheet (strategy_id-0) to add
# the "climate_change_factor_gnrl_hydropower_availability" variable. This
# pushes the existing rows 4:11 down to 5:12 (dimension grows to A1:AS12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(4).Insert()

# A4 / B4 - subsector / variable name
$ws.Cells.Item(4, 1).Value = "General"
$ws.Cells.Item(4, 2).Value = "climate_change_factor_gnrl_hydropower_availability"

# C4:G4 - blank (empty-string) metadata columns, matching every other data
# row on this sheet (normalize_group, trajgroup_no_vary_q, uniform_scaling_q,
# variable_trajectory_group, variable_trajectory_group_trajectory_type).
# Assigning "" directly clears the cell instead of leaving an empty text
# value, so use the text-prefix trick ("'") and then strip the resulting
# quote-prefix formatting it implies.
foreach ($col in 3..7) {
    $cell = $ws.Cells.Item(4, $col)
    $cell.Value = "'"
    $cell.ClearFormats()
}

# H4 / I4 - max_35 / min_35
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 0.5

# J4:AS4 - yearly values, all 1
for ($col = 10; $col -le 45; $col++) {
    $ws.Cells.Item(4, $col).Value = 1
}
